# IBP_TODO_AGILE.xlsx edit: add new tasks, update statuses, add autofilter
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 11: "use I cis target ..." task (pushes the old
#     "write a parsing script ..." row, and everything after it, down by one)
$ws.Rows(11).Insert()
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "use I cis target (find motifs, and parse so cbust can run on it)"

# --- Row 7: "run homer" task status is now finished
$ws.Range("D7").Value = "finished"

# --- Row 12 (was row 11 before insert): "write a parsing script ..." task
#     now has an owner + status, and its comment becomes the CRM-score note
$ws.Range("C12").Value = "Wim"
$ws.Range("D12").Value = "in progress"

# --- Row 13 (was row 12 before insert): rename to the final feature-matrix task
$ws.Range("B13").Value = "make the final feature matrix with all of our data combined, and maybe also send it to Daphne + Ibrahim"

# --- New row 19 appended at the end: Ibrahim's note about deep learning sequences
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Btw, Ibrahim mentioned that you can also start preparing the sequences for the deep learning models. Basically, you just have to make sure they all have the same length and are accompanied by a label"
$ws.Range("C19").Value = "Wim"
$ws.Range("B19").Font.Italic = $true
$ws.Range("B19").WrapText = $false

# --- Turn on the header autofilter, and register the (hidden) filter database name
$ws.Range("A1:F1").AutoFilter(1)
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$1")
$n.Visible = $false

# --- Move the active selection to B12 (matches where the user last worked)
$ws.Range("B12").Select()
